$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 48 / Row 49 swap (coin ranking changed: USDe now ranks above InjectiveProtocol) ---
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.78%  "

# --- Price / Volume(1h) refresh for all other rows ---
# Price (column D) values look numeric, so force text format (NumberFormat "@")
# before assigning them so Excel stores the exact original string (no float
# rounding / scientific notation), then restore the default "Normal" style so
# no stray formatting is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.199.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.75%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.260.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.41%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.75%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.260.25"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("E9").Value = "  -1.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.147"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000242"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.799.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.26%  "

$ws.Range("E16").Value = "  +0.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.260.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.275.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "470.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.728"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.94%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.103"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0725"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0398"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "425.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.055.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.119"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.259"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.33%  "

$ws.Range("E51").Value = "  -1.32%  "
